$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v1 = @'
Autopsies of deceased with a confirmed severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) infection can provide important insights into the novel disease and its course.
 Furthermore, autopsies are essential for the correct statistical recording of the coronavirus disease 2019 (COVID-19) deaths.
 In the northern German Federal State of Hamburg, all deaths of Hamburg citizens with ante- or postmortem PCR-confirmed SARS-CoV-2 infection have been autopsied since the outbreak of the pandemic in Germany.
 Our evaluation provides a systematic overview of the first 80 consecutive full autopsies.
 A proposal for the categorisation of deaths with SARS-CoV-2 infection is presented (category 1: definite COVID-19 death; category 2: probable COVID-19 death; category 3: possible COVID-19 death with an equal alternative cause of death; category 4: SARS-CoV-2 detection with cause of death not associated to COVID-19).
 In six cases, SARS-CoV-2 infection was diagnosed postmortem by a positive PCR test in a nasopharyngeal or lung tissue swab.
 In the other 74 cases, SARS-CoV-2 infection had already been known antemortem.
 The deceased were aged between 52 and 96 years (average 79.2 years, median 82.4 years).
 In the study cohort, 34 deceased were female (38%) and 46 male (62%).
 Overall, 38% of the deceased were overweight or obese.
 All deceased, except for two women, in whom no significant pre-existing conditions were found autoptically, had relevant comorbidities (in descending order of frequency): (1) diseases of the cardiovascular system, (2) lung diseases, (3) central nervous system diseases, (4) kidney diseases, and (5) diabetes mellitus.
 A total of 76 cases (95%) were classified as COVID-19 deaths, corresponding to categories 1–3. Four deaths (5%) were defined as non-COVID-19 deaths with virus-independent causes of death.
 In eight cases, pneumonia was combined with a fulminant pulmonary artery embolism.
 Peripheral pulmonary artery embolisms were found in nine other cases.
 Overall, deep vein thrombosis has been found in 40% of the cases.
 This study provides the largest overview of autopsies of SARS-CoV-2-infected patients presented so far.

'@
$ws.Range("D14").Value = $v1

$v2 = @'
[Carolin%Edler%NULL%1,  Ann Sophie%Schröder%NULL%0,  Martin%Aepfelbacher%NULL%0,  Antonia%Fitzek%NULL%1,  Axel%Heinemann%NULL%0,  Fabian%Heinrich%NULL%0,  Anke%Klein%NULL%1,  Felicia%Langenwalder%NULL%1,  Marc%Lütgehetmann%NULL%0,  Kira%Meißner%NULL%1,  Klaus%Püschel%NULL%0,  Julia%Schädler%NULL%1,  Stefan%Steurer%NULL%0,  Herbert%Mushumba%NULL%0,  Jan-Peter%Sperhake%sperhake@uke.de%0]
'@
$ws.Range("E14").Value = $v2

$v3 = @'
Forensic investigations generally contain extensive morphological examinations to accurately diagnose the cause of death.
 Thus, the appearance of a new disease often creates emerging challenges in morphological examinations due to the lack of available data from autopsy- or biopsy-based research.
 Since late December 2019, an outbreak of a novel seventh coronavirus disease has been reported in China caused by “severe acute respiratory syndrome coronavirus 2” (SARS-CoV-2).
 On March 11, 2020, the new clinical condition COVID-19 (Corona-Virus-Disease-19) was declared a pandemic by the World Health Organization (WHO).
 Patients with COVID-19 mainly have a mild disease course, but severe disease onset might result in death due to proceeded lung injury with massive alveolar damage and progressive respiratory failure.
 However, the detailed mechanisms that cause organ injury still remain unclear.
 We investigated the morphological findings of a COVID-19 patient who died during self-isolation.
 Pathologic examination revealed massive bilateral alveolar damage, indicating early-phase “acute respiratory distress syndrome” (ARDS).
 This case emphasizes the possibility of a rapid severe disease onset in previously mild clinical condition and highlights the necessity of a complete autopsy to gain a better understanding of the pathophysiological changes in SARS-CoV-2 infections.

'@
$ws.Range("D17").Value = $v3

$v4 = @'
[Christine%Suess%NULL%1,  Roland%Hausmann%roland.hausmann@kssg.ch%1]
'@
$ws.Range("E17").Value = $v4

$v5 = @'
In the setting of the coronavirus disease 2019 (COVID-19) pandemic, only few data regarding lung pathology induced by severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) is available, especially without medical intervention interfering with the natural evolution of the disease.
 We present here the first case of forensic autopsy of a COVID-19 fatality occurring in a young woman, in the community.
 Diagnosis was made at necropsy and lung histology showed diffuse alveolar damage, edema, and interstitial pneumonia with a geographically heterogeneous pattern, mostly affecting the central part of the lungs.
 This death related to COVID-19 pathology highlights the heterogeneity and severity of central lung lesions after natural evolution of the disease.

'@
$ws.Range("D18").Value = $v5

$v6 = @'
[Diego%Aguiar%diego.aguiar@hcuge.ch%1,  Johannes Alexander%Lobrinus%NULL%2,  Johannes Alexander%Lobrinus%NULL%0,  Manuel%Schibler%NULL%1,  Tony%Fracasso%NULL%1,  Christelle%Lardi%NULL%1]
'@
$ws.Range("E18").Value = $v6

$v7 = @'
A 75-year-old man presented to a French hospital with a 4-day fever after returning from a coronavirus disease-19 (COVID-19) cluster region.
 A reverse-transcription polymerase chain reaction test was positive for severe acute respiratory syndrome coronavirus-2 (SARS CoV-2) using a nasopharyngeal swab sample.
 After he returned home and a telephone follow-up, he was found deceased 9 days after first showing symptoms.
 Whole-body, non-enhanced, post-mortem computed tomography (PMCT) and a forensic autopsy were performed approximately 48 h after death, with sanitary precautions.
 The PMCT showed bilateral and diffuse crazy-paving lung opacities, with bilateral pleural effusions.
 Post-mortem virology studies detected the presence of SARS-CoV-2 (B.
1 lineage) in the nasopharynx, plasma, lung biopsies, pleural effusion and faeces confirming the persistence of viral ribonucleic acid 48 h after death.
 Microscopic examination showed that severe lung damage was responsible for his death.
 The main abnormality was diffuse alveolar damage, associated with different stages of inflammation and fibrosis.
 This case is one of the first to describe complete post-mortem data for a COVID-19 death and highlights the ability of PMCT to detect severe involvement of the lungs before autopsy in an apparently natural death.
 The present pathology results are concordant with previously reported findings and reinforce the disease pathogenesis hypothesis of combined viral replication with an inappropriate immune response.

'@
$ws.Range("D26").Value = $v7

$v8 = @'
[Mathilde%Ducloyer%Mathilde.ducloyer@chu-nantes.fr%1, Benjamin%Gaborit%NULL%2, Benjamin%Gaborit%NULL%0, Claire%Toquet%NULL%1, Louise%Castain%NULL%1, Antonin%Bal%NULL%1, Pierre Paul%Arrigoni%NULL%1, Raphaël%Lecomte%NULL%1, Renaud%Clement%NULL%1, Christine%Sagan%NULL%1]
'@
$ws.Range("E26").Value = $v8

$v9 = @'
We document the neuropathologic findings of a 73-year old man who died from acute cerebellar hemorrhage in the context of relatively mild SARS-CoV2 infection.
 The patient developed sudden onset of headache, nausea, and vomiting, immediately followed by loss of consciousness on the day of admission.
 Emergency medical services found him severely hypoxemic at home, and the patient suffered a cardiac arrest during transport to the emergency department.
 The emergency team achieved return of spontaneous circulation after over 17 min of resuscitation.
 A chest radiograph revealed hazy bilateral opacities; and real-time-PCR for SARS-CoV-2 on the nasopharyngeal swab was positive.
 Computed tomography of the head showed a large right cerebellar hemorrhage, with tonsillar herniation and intraventricular hemorrhage.
 One day after presentation, he was transitioned to comfort care and died shortly after palliative extubation.
 Autopsy performed 3 h after death showed cerebellar hemorrhage and acute infarcts in the dorsal pons and medulla.
 Remarkably, there were microglial nodules and neuronophagia bilaterally in the inferior olives and multifocally in the cerebellar dentate nuclei.
 This constellation of findings has not been reported thus far in the context of SARS-CoV-2 infection.

'@
$ws.Range("D35").Value = $v9

$v10 = @'
[Osama%Al-Dalahmah%oa2298@cumc.columbia.edu%1,  Kiran T.%Thakur%NULL%2,  Kiran T.%Thakur%NULL%0,  Anna S.%Nordvig%NULL%1,  Morgan L.%Prust%NULL%1,  William%Roth%NULL%1,  Angela%Lignelli%NULL%1,  Anne-Catrin%Uhlemann%NULL%1,  Emily Happy%Miller%NULL%1,  Shajo%Kunnath-Velayudhan%NULL%1,  Armando%Del Portillo%NULL%1,  Yang%Liu%NULL%1,  Gunnar%Hargus%NULL%1,  Andrew F.%Teich%NULL%1,  Richard A.%Hickman%NULL%1,  Kurenai%Tanji%NULL%1,  James E.%Goldman%NULL%1,  Phyllis L.%Faust%NULL%1,  Peter%Canoll%NULL%1]
'@
$ws.Range("E35").Value = $v10

$v11 = @'
This paper describes three autopsy cases with postmortem diagnosis of SARS-CoV-2 infection, with detailed macroscopic examination as well as advanced microscopic studies of organ tissues collected using hematoxylin-eosin stains and immunohistochemical markers.
 Two of the cases were admitted briefly in the County Clinical Emergency Hospital of Sibiu, and one was found deceased at his home address.
 All three autopsies were completed at the County morgue, in the COVID-19 restricted area, using complete protective equipment.
 The lungs of the patients seemed to be the center organ of invasion and pathogenesis of the novel coronavirus with diffuse areas of condensation, subpleural retraction zones but with different aspect of the classic bacterial bronchopneumonia.
 Microscopic evaluation revealed viral cytopathic effect of type II pneumocytes with a couple of cells that presented cytoplasmic and nuclear inclusions and who tend to form clusters mimicking multinucleated giant cells.
 Hyaline membranes and destruction of the alveolar wall as well as microthrombi formation within the small blood vessels were constantly found in almost all our three cases.
 The spleen had sustained white pulp atrophy with absence of lymphoid follicles.
 There were no microscopic signs of viral infection on the myocardium or the other organs.

'@
$ws.Range("D36").Value = $v11

$v12 = @'
[George-Călin%Oprinca%george_oprinca@yahoo.com%1,  Lilioara-Alexandra%Muja%NULL%2,  Lilioara-Alexandra%Muja%NULL%0]
'@
$ws.Range("E36").Value = $v12

$v13 = @'
Background
Post-mortem studies can provide important information for understanding new diseases and small autopsy case series have already reported different findings in COVID-19 patients.
Methods
We evaluated whether some specific post-mortem features are observed in these patients and if these changes are related to the presence of the virus in different organs.
 Complete macroscopic and microscopic autopsies were performed on different organs in 17 COVID-19 non-survivors.
 Presence of SARS-CoV-2 was evaluated with immunohistochemistry (IHC) in lung samples and with real-time reverse-transcription polymerase chain reaction (RT-PCR) test in the lung and other organs.
Results
Pulmonary findings revealed early-stage diffuse alveolar damage (DAD) in 15 out of 17 patients and microthrombi in small lung arteries in 11 patients.
 Late-stage DAD, atypical pneumocytes, and/or acute pneumonia were also observed.
 Four lung infarcts, two acute myocardial infarctions, and one ischemic enteritis were observed.
 There was no evidence of myocarditis, hepatitis, or encephalitis.
 Kidney evaluation revealed the presence of hemosiderin in tubules or pigmented casts in most patients.
 Spongiosis and vascular congestion were the most frequently encountered brain lesions.
 No specific SARS-CoV-2 lesions were observed in any organ.
 IHC revealed positive cells with a heterogeneous distribution in the lungs of 11 of the 17 (65%) patients; RT-PCR yielded a wide distribution of SARS-CoV-2 in different tissues, with 8 patients showing viral presence in all tested organs (i.
e.
, lung, heart, spleen, liver, colon, kidney, and brain).
Conclusions
In conclusion, autopsies revealed a great heterogeneity of COVID-19-associated organ injury and the remarkable absence of any specific viral lesions, even when RT-PCR identified the presence of the virus in many organs.

'@
$ws.Range("D31").Value = $v13

$v14 = @'
[Myriam%Remmelink%NULL%1,     Ricardo%De Mendonça%NULL%1,     Nicky%D’Haene%NULL%1,     Sarah%De Clercq%NULL%1,     Camille%Verocq%NULL%1,     Laetitia%Lebrun%NULL%1,     Philomène%Lavis%NULL%1,     Marie-Lucie%Racu%NULL%1,     Anne-Laure%Trépant%NULL%1,     Calliope%Maris%NULL%1,     Sandrine%Rorive%NULL%1,     Jean-Christophe%Goffard%NULL%1,     Olivier%De Witte%NULL%1,     Lorenzo%Peluso%NULL%1,     Jean-Louis%Vincent%NULL%1,     Christine%Decaestecker%NULL%1,     Fabio Silvio%Taccone%NULL%1,     Isabelle%Salmon%Isabelle.Salmon@erasme.ulb.ac.be%1]
'@
$ws.Range("E31").Value = $v14

